$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-533) holds a "Förändrad" date that was bulk-updated
# from serial 45182 (2023-09-13) to 45184 (2023-09-15).
$ws.Range("C2:C533").Value = 45184
